# Prefix each protocol sheet's Step/command names in column A with the
# sheet (protocol) name, e.g. "Step4 Takeaway" -> "discount2 Step4 Takeaway".
# This mirrors the commit: "fix: unique command names in XLSX - prefix
# protocol name to each step".

$wb = $excel.ActiveWorkbook

$sheetNames = @(
    "price1", "price2",
    "discount1", "discount2",
    "free1", "free2",
    "nomoney1", "nomoney2",
    "noppv1", "noppv2",
    "card1", "card2",
    "nosex1", "nosex2",
    "offtopic1", "offtopic2",
    "real1", "real2",
    "voice1", "voice2",
    "customyes1", "customyes2",
    "customno1", "customno2",
    "done1", "done2",
    "cumcontrol", "dickpic", "boosters"
)

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    # Find the last used row in column A (data starts at row 2; row 1 is the header "Name").
    $lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row

    for ($r = 2; $r -le $lastRow; $r++) {
        $cell = $ws.Cells.Item($r, 1)
        $current = $cell.Value2
        if ($current -ne $null -and $current -ne "") {
            $cell.Value = "$sheetName $current"
        }
    }
}
